$d = $word.ActiveDocument

# --- Step 1: the "Cuarta" paragraph currently holds its text split across two
# runs ("Cuarta" + " página web de la empresa"). Collapse it back down to a
# single run by re-writing the matched span through Find/Replace (Word merges
# the replaced text into one run using the first run's formatting).
$rng = $d.Content
$rng.Find.Execute("Cuarta página web de la empresa", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Cuarta página web de la empresa", 2) | Out-Null

# --- Step 2: append a brand-new "Quinto página web de la empresa" paragraph
# right after it, reproducing the same two-run split ("Quinto" + " página web
# de la empresa") the "Cuarta" paragraph originally had, followed by one more
# empty paragraph. InsertXML lets us author the exact run layout instead of
# relying on Word's auto-merge-on-type behaviour.
$cuarta = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "Cuarta página web de la empresa") {
        $cuarta = $d.Paragraphs($i)
        break
    }
}
$tail = $d.Range($cuarta.Range.End, $cuarta.Range.End)

$newXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:r><w:t>Quinto</w:t></w:r><w:r><w:t xml:space="preserve"> página web de la empresa</w:t></w:r></w:p><w:p/></w:body>' + `
    '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$tail.InsertXML($newXml) | Out-Null
